$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (B, C, D, E, G) - column F (Win) is unchanged
$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    3 = @(0.3464964993005633, 9.226618575922256, 3.082599426703578, 6.48142807727062, 19.13714257919702)
    4 = @(0.02258322285507441, 0.3375848360084654, 3.082599426703578, 6.48142807727062, 9.924195562837738)
    5 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    6 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
